# Append 20 new departure rows (206-225) to the "Main Data" sheet, mirroring
# the existing table layout: A=NUMBER, B=DATE, C=TIME, D=FLIGHT, E=TO,
# F=SHORT, G=AIRLINE, H=MODEL, I=AIRCFAT ID, J=STATUS, K=(blank),
# L=DIFFERENCE, M=(blank).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A206").Value2 = 205.0
$ws.Range("B206").Value2 = 'Friday, Jan 13'
$ws.Range("C206").Value2 = '4:00 AM'
$ws.Range("D206").Value2 = 'CGF5919'
$ws.Range("E206").Value2 = 'Leipzig'
$ws.Range("F206").Value2 = '(LEJ)'
$ws.Range("G206").Value2 = 'Cargo Air '
$ws.Range("H206").Value2 = 'B733'
$ws.Range("I206").Value2 = '(LZ-CGP)'
$ws.Range("J206").Value2 = '3:45 AM'
$ws.Range("L206").Value2 = '0 hours, -15 minutes'

$ws.Range("A207").Value2 = 206.0
$ws.Range("B207").Value2 = 'Friday, Jan 13'
$ws.Range("C207").Value2 = '5:10 AM'
$ws.Range("D207").Value2 = 'BO625'
$ws.Range("E207").Value2 = 'Madrid'
$ws.Range("F207").Value2 = '(MAD)'
$ws.Range("G207").Value2 = 'Bluebird Nordic '
$ws.Range("H207").Value2 = 'B734'
$ws.Range("I207").Value2 = '(TF-BBN)'
$ws.Range("J207").Value2 = '5:14 AM'
$ws.Range("L207").Value2 = '0 hours, 4 minutes'

$ws.Range("A208").Value2 = 207.0
$ws.Range("B208").Value2 = 'Friday, Jan 13'
$ws.Range("C208").Value2 = '5:40 AM'
$ws.Range("D208").Value2 = 'LO3880'
$ws.Range("E208").Value2 = 'Warsaw'
$ws.Range("F208").Value2 = '(WAW)'
$ws.Range("G208").Value2 = 'LOT '
$ws.Range("H208").Value2 = 'E75S'
$ws.Range("I208").Value2 = '(SP-LIK)'
$ws.Range("J208").Value2 = '5:47 AM'
$ws.Range("L208").Value2 = '0 hours, 7 minutes'

$ws.Range("A209").Value2 = 208.0
$ws.Range("B209").Value2 = 'Friday, Jan 13'
$ws.Range("C209").Value2 = '6:15 AM'
$ws.Range("D209").Value2 = 'W61001'
$ws.Range("E209").Value2 = 'London'
$ws.Range("F209").Value2 = '(LTN)'
$ws.Range("G209").Value2 = 'Wizz Air '
$ws.Range("H209").Value2 = 'A321'
$ws.Range("I209").Value2 = '(HA-LXE)'
$ws.Range("J209").Value2 = '6:21 AM'
$ws.Range("L209").Value2 = '0 hours, 6 minutes'

$ws.Range("A210").Value2 = 209.0
$ws.Range("B210").Value2 = 'Friday, Jan 13'
$ws.Range("C210").Value2 = '6:20 AM'
$ws.Range("D210").Value2 = 'FR6389'
$ws.Range("E210").Value2 = 'Dortmund'
$ws.Range("F210").Value2 = '(DTM)'
$ws.Range("G210").Value2 = 'Ryanair '
$ws.Range("H210").Value2 = 'B738'
$ws.Range("I210").Value2 = '(SP-RKB)'
$ws.Range("J210").Value2 = '6:29 AM'
$ws.Range("L210").Value2 = '0 hours, 9 minutes'

$ws.Range("A211").Value2 = 210.0
$ws.Range("B211").Value2 = 'Friday, Jan 13'
$ws.Range("C211").Value2 = '6:20 AM'
$ws.Range("D211").Value2 = 'W61043'
$ws.Range("E211").Value2 = 'Catania'
$ws.Range("F211").Value2 = '(CTA)'
$ws.Range("G211").Value2 = 'Wizz Air '
$ws.Range("H211").Value2 = 'A321'
$ws.Range("I211").Value2 = '(HA-LXP)'
$ws.Range("J211").Value2 = '6:41 AM'
$ws.Range("L211").Value2 = '0 hours, 21 minutes'

$ws.Range("A212").Value2 = 211.0
$ws.Range("B212").Value2 = 'Friday, Jan 13'
$ws.Range("C212").Value2 = '6:30 AM'
$ws.Range("D212").Value2 = 'W61091'
$ws.Range("E212").Value2 = 'Dortmund'
$ws.Range("F212").Value2 = '(DTM)'
$ws.Range("G212").Value2 = 'Wizz Air '
$ws.Range("H212").Value2 = 'A320'
$ws.Range("I212").Value2 = '(HA-LYK)'
$ws.Range("J212").Value2 = '6:37 AM'
$ws.Range("L212").Value2 = '0 hours, 7 minutes'

$ws.Range("A213").Value2 = 212.0
$ws.Range("B213").Value2 = 'Friday, Jan 13'
$ws.Range("C213").Value2 = '6:35 AM'
$ws.Range("D213").Value2 = '3Z7606'
$ws.Range("E213").Value2 = 'Salalah'
$ws.Range("F213").Value2 = '(SLL)'
$ws.Range("G213").Value2 = 'Smartwings '
$ws.Range("H213").Value2 = 'B38M'
$ws.Range("I213").Value2 = '(OK-SWC)'
$ws.Range("J213").Value2 = '7:00 AM'
$ws.Range("L213").Value2 = '0 hours, 25 minutes'

$ws.Range("A214").Value2 = 213.0
$ws.Range("B214").Value2 = 'Friday, Jan 13'
$ws.Range("C214").Value2 = '6:40 AM'
$ws.Range("D214").Value2 = 'W61163'
$ws.Range("E214").Value2 = 'Malta'
$ws.Range("F214").Value2 = '(MLA)'
$ws.Range("G214").Value2 = 'Wizz Air '
$ws.Range("H214").Value2 = 'A321'
$ws.Range("I214").Value2 = '(HA-LTC)'
$ws.Range("J214").Value2 = '6:52 AM'
$ws.Range("L214").Value2 = '0 hours, 12 minutes'

$ws.Range("A215").Value2 = 214.0
$ws.Range("B215").Value2 = 'Friday, Jan 13'
$ws.Range("C215").Value2 = '6:45 AM'
$ws.Range("D215").Value2 = 'LH1363'
$ws.Range("E215").Value2 = 'Frankfurt'
$ws.Range("F215").Value2 = '(FRA)'
$ws.Range("G215").Value2 = 'Lufthansa '
$ws.Range("H215").Value2 = 'CRJ9'
$ws.Range("I215").Value2 = '(D-ACNK)'
$ws.Range("J215").Value2 = '6:44 AM'
$ws.Range("L215").Value2 = '0 hours, -1 minutes'

$ws.Range("A216").Value2 = 215.0
$ws.Range("B216").Value2 = 'Friday, Jan 13'
$ws.Range("C216").Value2 = '7:10 AM'
$ws.Range("D216").Value2 = 'W61071'
$ws.Range("E216").Value2 = 'Eindhoven'
$ws.Range("F216").Value2 = '(EIN)'
$ws.Range("G216").Value2 = 'Wizz Air '
$ws.Range("H216").Value2 = 'A321'
$ws.Range("I216").Value2 = '(HA-LXD)'
$ws.Range("J216").Value2 = '7:17 AM'
$ws.Range("L216").Value2 = '0 hours, 7 minutes'

$ws.Range("A217").Value2 = 216.0
$ws.Range("B217").Value2 = 'Friday, Jan 13'
$ws.Range("C217").Value2 = '8:30 AM'
$ws.Range("D217").Value2 = '3Z7632'
$ws.Range("E217").Value2 = 'Hurghada'
$ws.Range("F217").Value2 = '(HRG)'
$ws.Range("G217").Value2 = 'Smartwings '
$ws.Range("H217").Value2 = 'B738'
$ws.Range("I217").Value2 = '(HA-LKG)'
$ws.Range("J217").Value2 = '8:37 AM'
$ws.Range("L217").Value2 = '0 hours, 7 minutes'

$ws.Range("A218").Value2 = 217.0
$ws.Range("B218").Value2 = 'Friday, Jan 13'
$ws.Range("C218").Value2 = '9:35 AM'
$ws.Range("D218").Value2 = 'BO951'
$ws.Range("E218").Value2 = 'Paris'
$ws.Range("F218").Value2 = '(CDG)'
$ws.Range("G218").Value2 = 'Bluebird Nordic '
$ws.Range("H218").Value2 = 'B734'
$ws.Range("I218").Value2 = '(TF-BBJ)'
$ws.Range("J218").Value2 = '9:26 AM'
$ws.Range("L218").Value2 = '0 hours, -9 minutes'

$ws.Range("A219").Value2 = 218.0
$ws.Range("B219").Value2 = 'Friday, Jan 13'
$ws.Range("C219").Value2 = '9:55 AM'
$ws.Range("D219").Value2 = 'E47011'
$ws.Range("E219").Value2 = 'Dubai'
$ws.Range("F219").Value2 = '(DWC)'
$ws.Range("G219").Value2 = 'Enter Air '
$ws.Range("H219").Value2 = 'B738'
$ws.Range("I219").Value2 = '(SP-ENW)'
$ws.Range("J219").Value2 = '10:09 AM'
$ws.Range("L219").Value2 = '0 hours, 14 minutes'

$ws.Range("A220").Value2 = 219.0
$ws.Range("B220").Value2 = 'Friday, Jan 13'
$ws.Range("C220").Value2 = '10:00 AM'
$ws.Range("D220").Value2 = 'RR7989'
$ws.Range("E220").Value2 = 'Gran Canaria'
$ws.Range("F220").Value2 = '(LPA)'
$ws.Range("G220").Value2 = 'Ryanair '
$ws.Range("H220").Value2 = 'B738'
$ws.Range("I220").Value2 = '(SP-RSN)'
$ws.Range("J220").Value2 = '9:54 AM'
$ws.Range("L220").Value2 = '0 hours, -6 minutes'

$ws.Range("A221").Value2 = 220.0
$ws.Range("B221").Value2 = 'Friday, Jan 13'
$ws.Range("C221").Value2 = '10:25 AM'
$ws.Range("D221").Value2 = 'FR2472'
$ws.Range("E221").Value2 = 'London'
$ws.Range("F221").Value2 = '(STN)'
$ws.Range("G221").Value2 = 'Ryanair '
$ws.Range("H221").Value2 = 'B738'
$ws.Range("I221").Value2 = '(SP-RKB)'
$ws.Range("J221").Value2 = '11:04 AM'
$ws.Range("L221").Value2 = '0 hours, 39 minutes'

$ws.Range("A222").Value2 = 221.0
$ws.Range("B222").Value2 = 'Friday, Jan 13'
$ws.Range("C222").Value2 = '11:45 AM'
$ws.Range("D222").Value2 = 'LO3882'
$ws.Range("E222").Value2 = 'Warsaw'
$ws.Range("F222").Value2 = '(WAW)'
$ws.Range("G222").Value2 = 'LOT '
$ws.Range("H222").Value2 = 'E75S'
$ws.Range("I222").Value2 = '(SP-LIA)'
$ws.Range("J222").Value2 = '11:48 AM'
$ws.Range("L222").Value2 = '0 hours, 3 minutes'

$ws.Range("A223").Value2 = 222.0
$ws.Range("B223").Value2 = 'Friday, Jan 13'
$ws.Range("C223").Value2 = '12:40 PM'
$ws.Range("D223").Value2 = 'FR6892'
$ws.Range("E223").Value2 = 'Dortmund'
$ws.Range("F223").Value2 = '(DTM)'
$ws.Range("G223").Value2 = 'Ryanair '
$ws.Range("H223").Value2 = 'B738'
$ws.Range("I223").Value2 = '(SP-RSB)'
$ws.Range("J223").Value2 = '12:51 PM'
$ws.Range("L223").Value2 = '0 hours, 11 minutes'

$ws.Range("A224").Value2 = 223.0
$ws.Range("B224").Value2 = 'Friday, Jan 13'
$ws.Range("C224").Value2 = '12:45 PM'
$ws.Range("D224").Value2 = 'W61275'
$ws.Range("E224").Value2 = 'Abu Dhabi'
$ws.Range("F224").Value2 = '(AUH)'
$ws.Range("G224").Value2 = 'Wizz Air '
$ws.Range("H224").Value2 = 'A21N'
$ws.Range("I224").Value2 = '(HA-LZJ)'
$ws.Range("J224").Value2 = '12:48 PM'
$ws.Range("L224").Value2 = '0 hours, 3 minutes'

$ws.Range("A225").Value2 = 224.0
$ws.Range("B225").Value2 = 'Friday, Jan 13'
$ws.Range("C225").Value2 = '2:15 PM'
$ws.Range("D225").Value2 = 'LH1357'
$ws.Range("E225").Value2 = 'Frankfurt'
$ws.Range("F225").Value2 = '(FRA)'
$ws.Range("G225").Value2 = 'Lufthansa '
$ws.Range("H225").Value2 = 'CRJ9'
$ws.Range("I225").Value2 = '(D-ACKI)'
$ws.Range("J225").Value2 = '2:15 PM'
$ws.Range("L225").Value2 = '0 hours, 0 minutes'

# K and M are blank spacer columns in every existing row; copy the blank,
# default-styled cells from row 2 so the new rows keep the same shape.
$ws.Range("K2").Copy($ws.Range("K206:K225"))
$ws.Range("M2").Copy($ws.Range("M206:M225"))

Write-Output "Rows 206-225 added successfully"
